$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 687, pushing existing rows 687.. down by one
$ws.Rows.Item(687).Insert()

# Populate the newly inserted row with the new data record
$ws.Cells.Item(687, 1).Value = 10
$ws.Cells.Item(687, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(687, 3).Value = "La Araucanía"
$ws.Cells.Item(687, 4).Value = 45154
$ws.Cells.Item(687, 5).Value = 9
$ws.Cells.Item(687, 6).Value = "Fruta"
$ws.Cells.Item(687, 7).Value = 100108
$ws.Cells.Item(687, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(687, 9).Value = 100108005
$ws.Cells.Item(687, 10).Value = "Piña"
$ws.Cells.Item(687, 11).Value = "Caramelo"
$ws.Cells.Item(687, 12).Value = "Primera"
$ws.Cells.Item(687, 13).Value = 55
$ws.Cells.Item(687, 14).Value = 26000
$ws.Cells.Item(687, 15).Value = 26000
$ws.Cells.Item(687, 16).Value = 26000
$ws.Cells.Item(687, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(687, 18).Value = "Ecuador"
$ws.Cells.Item(687, 19).Value = 2167
$ws.Cells.Item(687, 20).Value = 12
